# Wiring.pptx — "new head and tail"
#
# 1. Refresh the cached "datetimeFigureOut" date field (14.05.2022 -> 20.05.2022)
#    on the slide master and on every slide layout's Date Placeholder.
# 2. Re-crop / re-position the first picture ("Picture 4") on slide 4 to show a
#    new head-and-tail crop of the source image.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text (slide master + all slide layouts)
# ---------------------------------------------------------------------------
$newDate = "20.05.2022"
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lay = $layouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 4 - "Picture 4": new crop + new position/size
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$pic = $slide4.Shapes.Item(1)

# Crop (points, relative to the native image size) -> a:srcRect l="16293" t="23938"
$pic.PictureFormat.CropLeft = 144.926235
$pic.PictureFormat.CropTop = 158.888475

# New position/size (points) -> a:off x="444500" y="2043472" a:ext cx="6557397" cy="4446228"
$pic.Left = 35.0
$pic.Top = 160.90331268661384
$pic.Width = 516.3305
$pic.Height = 350.0967

Write-Output "done"
